$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 0.2679253333333333
$ws.Cells.Item(2, 8).Value2 = 0.8037759999999999
$ws.Cells.Item(2, 9).Value2 = 0.1226600350746756
$ws.Cells.Item(2, 10).Value2 = 0.1226600350746756
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 0.4702473333333333
$ws.Cells.Item(2, 14).Value2 = 1.410742
$ws.Cells.Item(2, 15).Value2 = 0.009034922268422819
$ws.Cells.Item(2, 16).Value2 = 0.009034922268422819
$ws.Cells.Item(2, 17).Value2 = 0.1259911735324444
$ws.Cells.Item(2, 18).Value2 = 1.133920561792
$ws.Cells.Item(2, 19).Value2 = 0.001108223882341711
$ws.Cells.Item(2, 20).Value2 = 0.001108223882341711
# Row 3
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 0.2679253333333333
$ws.Cells.Item(3, 8).Value2 = 0.8037759999999999
$ws.Cells.Item(3, 9).Value2 = 0.1226600350746756
$ws.Cells.Item(3, 10).Value2 = 0.1226600350746756
$ws.Cells.Item(3, 14).Value2 = 0.9584440000000001
$ws.Cells.Item(3, 15).Value2 = 0.006138235792679485
$ws.Cells.Item(3, 16).Value2 = 0.006138235792679485
$ws.Cells.Item(3, 17).Value2 = 0.08559714272711109
$ws.Cells.Item(3, 18).Value2 = 0.770374284544
$ws.Cells.Item(3, 19).Value2 = 0.0007529162176266947
$ws.Cells.Item(3, 20).Value2 = 0.0007529162176266949
# Row 4
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 0.2679253333333333
$ws.Cells.Item(4, 8).Value2 = 0.8037759999999999
$ws.Cells.Item(4, 9).Value2 = 0.1226600350746756
$ws.Cells.Item(4, 10).Value2 = 0.1226600350746756
$ws.Cells.Item(4, 13).Value2 = 1.047307
$ws.Cells.Item(4, 14).Value2 = 3.141921
$ws.Cells.Item(4, 15).Value2 = 0.02012204358311108
$ws.Cells.Item(4, 16).Value2 = 0.02012204358311108
$ws.Cells.Item(4, 17).Value2 = 0.2806000770773333
$ws.Cells.Item(4, 18).Value2 = 2.525400693696
$ws.Cells.Item(4, 19).Value2 = 0.002468170571678556
$ws.Cells.Item(4, 20).Value2 = 0.002468170571678557
# Row 5
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 0.2679253333333333
$ws.Cells.Item(5, 8).Value2 = 0.8037759999999999
$ws.Cells.Item(5, 9).Value2 = 0.1226600350746756
$ws.Cells.Item(5, 10).Value2 = 0.1226600350746756
$ws.Cells.Item(5, 13).Value2 = 50.21070966666667
$ws.Cells.Item(5, 14).Value2 = 150.632129
$ws.Cells.Item(5, 15).Value2 = 0.9647047983557866
$ws.Cells.Item(5, 16).Value2 = 0.9647047983557866
$ws.Cells.Item(5, 17).Value2 = 13.45272112434489
$ws.Cells.Item(5, 18).Value2 = 121.074490119104
$ws.Cells.Item(5, 19).Value2 = 0.1183307244030286
$ws.Cells.Item(5, 20).Value2 = 0.1183307244030287
# Row 6
$ws.Cells.Item(6, 9).Value2 = 0.327101565785771
$ws.Cells.Item(6, 10).Value2 = 0.327101565785771
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 0.4702473333333333
$ws.Cells.Item(6, 14).Value2 = 1.410742
$ws.Cells.Item(6, 15).Value2 = 0.009034922268422819
$ws.Cells.Item(6, 16).Value2 = 0.009034922268422819
$ws.Cells.Item(6, 17).Value2 = 0.3359848227057777
$ws.Cells.Item(6, 18).Value2 = 3.023863404352
$ws.Cells.Item(6, 19).Value2 = 0.002955337220753834
$ws.Cells.Item(6, 20).Value2 = 0.002955337220753835
# Row 7
$ws.Cells.Item(7, 9).Value2 = 0.327101565785771
$ws.Cells.Item(7, 10).Value2 = 0.327101565785771
$ws.Cells.Item(7, 14).Value2 = 0.9584440000000001
$ws.Cells.Item(7, 15).Value2 = 0.006138235792679485
$ws.Cells.Item(7, 16).Value2 = 0.006138235792679485
$ws.Cells.Item(7, 19).Value2 = 0.002007826538947723
$ws.Cells.Item(7, 20).Value2 = 0.002007826538947723
# Row 8
$ws.Cells.Item(8, 9).Value2 = 0.327101565785771
$ws.Cells.Item(8, 10).Value2 = 0.327101565785771
$ws.Cells.Item(8, 13).Value2 = 1.047307
$ws.Cells.Item(8, 14).Value2 = 3.141921
$ws.Cells.Item(8, 15).Value2 = 0.02012204358311108
$ws.Cells.Item(8, 16).Value2 = 0.02012204358311108
$ws.Cells.Item(8, 17).Value2 = 0.7482854909973333
$ws.Cells.Item(8, 18).Value2 = 6.734569418976
$ws.Cells.Item(8, 19).Value2 = 0.006581951962845161
$ws.Cells.Item(8, 20).Value2 = 0.006581951962845162
# Row 9
$ws.Cells.Item(9, 9).Value2 = 0.327101565785771
$ws.Cells.Item(9, 10).Value2 = 0.327101565785771
$ws.Cells.Item(9, 13).Value2 = 50.21070966666667
$ws.Cells.Item(9, 14).Value2 = 150.632129
$ws.Cells.Item(9, 15).Value2 = 0.9647047983557866
$ws.Cells.Item(9, 16).Value2 = 0.9647047983557866
$ws.Cells.Item(9, 17).Value2 = 35.87481563309156
$ws.Cells.Item(9, 18).Value2 = 322.8733406978241
$ws.Cells.Item(9, 19).Value2 = 0.3155564500632242
$ws.Cells.Item(9, 20).Value2 = 0.3155564500632243
# Row 10
$ws.Cells.Item(10, 7).Value2 = 1.145196333333333
$ws.Cells.Item(10, 8).Value2 = 3.435589
$ws.Cells.Item(10, 9).Value2 = 0.5242871984759059
$ws.Cells.Item(10, 10).Value2 = 0.5242871984759059
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 12).Value2 = 1
$ws.Cells.Item(10, 13).Value2 = 0.4702473333333333
$ws.Cells.Item(10, 14).Value2 = 1.410742
$ws.Cells.Item(10, 15).Value2 = 0.009034922268422819
$ws.Cells.Item(10, 16).Value2 = 0.009034922268422819
$ws.Cells.Item(10, 17).Value2 = 0.5385255218931111
$ws.Cells.Item(10, 18).Value2 = 4.846729697038
$ws.Cells.Item(10, 19).Value2 = 0.004736894084558976
$ws.Cells.Item(10, 20).Value2 = 0.004736894084558976
# Row 11
$ws.Cells.Item(11, 7).Value2 = 1.145196333333333
$ws.Cells.Item(11, 8).Value2 = 3.435589
$ws.Cells.Item(11, 9).Value2 = 0.5242871984759059
$ws.Cells.Item(11, 10).Value2 = 0.5242871984759059
$ws.Cells.Item(11, 14).Value2 = 0.9584440000000001
$ws.Cells.Item(11, 15).Value2 = 0.006138235792679485
$ws.Cells.Item(11, 16).Value2 = 0.006138235792679485
$ws.Cells.Item(11, 17).Value2 = 0.3658688515017778
$ws.Cells.Item(11, 18).Value2 = 3.292819663516001
$ws.Cells.Item(11, 19).Value2 = 0.003218198447328459
$ws.Cells.Item(11, 20).Value2 = 0.003218198447328459
# Row 12
$ws.Cells.Item(12, 7).Value2 = 1.145196333333333
$ws.Cells.Item(12, 8).Value2 = 3.435589
$ws.Cells.Item(12, 9).Value2 = 0.5242871984759059
$ws.Cells.Item(12, 10).Value2 = 0.5242871984759059
$ws.Cells.Item(12, 13).Value2 = 1.047307
$ws.Cells.Item(12, 14).Value2 = 3.141921
$ws.Cells.Item(12, 15).Value2 = 0.02012204358311108
$ws.Cells.Item(12, 16).Value2 = 0.02012204358311108
$ws.Cells.Item(12, 17).Value2 = 1.199372136274333
$ws.Cells.Item(12, 18).Value2 = 10.794349226469
$ws.Cells.Item(12, 19).Value2 = 0.01054972985779939
$ws.Cells.Item(12, 20).Value2 = 0.01054972985779939
# Row 13
$ws.Cells.Item(13, 7).Value2 = 1.145196333333333
$ws.Cells.Item(13, 8).Value2 = 3.435589
$ws.Cells.Item(13, 9).Value2 = 0.5242871984759059
$ws.Cells.Item(13, 10).Value2 = 0.5242871984759059
$ws.Cells.Item(13, 13).Value2 = 50.21070966666667
$ws.Cells.Item(13, 14).Value2 = 150.632129
$ws.Cells.Item(13, 15).Value2 = 0.9647047983557866
$ws.Cells.Item(13, 16).Value2 = 0.9647047983557866
$ws.Cells.Item(13, 17).Value2 = 57.50112060433123
$ws.Cells.Item(13, 18).Value2 = 517.5100854389812
$ws.Cells.Item(13, 19).Value2 = 0.505782376086219
$ws.Cells.Item(13, 20).Value2 = 0.505782376086219
# Row 14
$ws.Cells.Item(14, 5).Value2 = 2
$ws.Cells.Item(14, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(14, 7).Value2 = 0.05668500000000001
$ws.Cells.Item(14, 8).Value2 = 0.170055
$ws.Cells.Item(14, 9).Value2 = 0.02595120066364754
$ws.Cells.Item(14, 10).Value2 = 0.02595120066364754
$ws.Cells.Item(14, 11).Value2 = 3
$ws.Cells.Item(14, 12).Value2 = 1
$ws.Cells.Item(14, 13).Value2 = 0.4702473333333333
$ws.Cells.Item(14, 14).Value2 = 1.410742
$ws.Cells.Item(14, 15).Value2 = 0.009034922268422819
$ws.Cells.Item(14, 16).Value2 = 0.009034922268422819
$ws.Cells.Item(14, 17).Value2 = 0.02665597009
$ws.Cells.Item(14, 18).Value2 = 0.23990373081
$ws.Cells.Item(14, 19).Value2 = 0.0002344670807682982
$ws.Cells.Item(14, 20).Value2 = 0.0002344670807682982
# Row 15
$ws.Cells.Item(15, 5).Value2 = 2
$ws.Cells.Item(15, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(15, 7).Value2 = 0.05668500000000001
$ws.Cells.Item(15, 8).Value2 = 0.170055
$ws.Cells.Item(15, 9).Value2 = 0.02595120066364754
$ws.Cells.Item(15, 10).Value2 = 0.02595120066364754
$ws.Cells.Item(15, 14).Value2 = 0.9584440000000001
$ws.Cells.Item(15, 15).Value2 = 0.006138235792679485
$ws.Cells.Item(15, 16).Value2 = 0.006138235792679485
$ws.Cells.Item(15, 17).Value2 = 0.01810979938
$ws.Cells.Item(15, 18).Value2 = 0.16298819442
$ws.Cells.Item(15, 19).Value2 = 0.0001592945887766089
$ws.Cells.Item(15, 20).Value2 = 0.0001592945887766089
# Row 16
$ws.Cells.Item(16, 5).Value2 = 2
$ws.Cells.Item(16, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(16, 7).Value2 = 0.05668500000000001
$ws.Cells.Item(16, 8).Value2 = 0.170055
$ws.Cells.Item(16, 9).Value2 = 0.02595120066364754
$ws.Cells.Item(16, 10).Value2 = 0.02595120066364754
$ws.Cells.Item(16, 13).Value2 = 1.047307
$ws.Cells.Item(16, 14).Value2 = 3.141921
$ws.Cells.Item(16, 15).Value2 = 0.02012204358311108
$ws.Cells.Item(16, 16).Value2 = 0.02012204358311108
$ws.Cells.Item(16, 17).Value2 = 0.059366597295
$ws.Cells.Item(16, 18).Value2 = 0.5342993756550001
$ws.Cells.Item(16, 19).Value2 = 0.000522191190787977
$ws.Cells.Item(16, 20).Value2 = 0.000522191190787977
# Row 17
$ws.Cells.Item(17, 5).Value2 = 2
$ws.Cells.Item(17, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(17, 7).Value2 = 0.05668500000000001
$ws.Cells.Item(17, 8).Value2 = 0.170055
$ws.Cells.Item(17, 9).Value2 = 0.02595120066364754
$ws.Cells.Item(17, 10).Value2 = 0.02595120066364754
$ws.Cells.Item(17, 13).Value2 = 50.21070966666667
$ws.Cells.Item(17, 14).Value2 = 150.632129
$ws.Cells.Item(17, 15).Value2 = 0.9647047983557866
$ws.Cells.Item(17, 16).Value2 = 0.9647047983557866
$ws.Cells.Item(17, 17).Value2 = 2.846194077455001
$ws.Cells.Item(17, 18).Value2 = 25.61574669709501
$ws.Cells.Item(17, 19).Value2 = 0.02503524780331465
$ws.Cells.Item(17, 20).Value2 = 0.02503524780331465
